$p = $ppt.ActivePresentation
$master = $p.Designs.Item(1).SlideMaster
$sh = $master.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange
Write-Host "Text: [$($tr.Text)] Length=$($tr.Length)"
$c = $tr.Characters(1,2)
Write-Host "Chars(1,2): [$($c.Text)]"
$c.Text = "19"
Write-Host "Text after small edit: [$($tr.Text)]"
